# Scheduled market-price refresh: update currentAveragePrice / Leve price / profit
# columns (H:N) across the per-job sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1562.5319
$ws.Range("I15").Value = 1562.5319
$ws.Range("K15").Value = 4687.5957
$ws.Range("M15").Value = -4518.5957
# Row 113
$ws.Range("H113").Value = 4420.028
$ws.Range("I113").Value = 4331.875
$ws.Range("J113").Value = 4596.3335
$ws.Range("K113").Value = 4331.875
$ws.Range("L113").Value = 4596.3335
$ws.Range("M113").Value = -1077.875
$ws.Range("N113").Value = -11104.3335
# Row 135
$ws.Range("H135").Value = 23810354
$ws.Range("I135").Value = 671.84375
$ws.Range("J135").Value = 100001340
$ws.Range("K135").Value = 6046.59375
$ws.Range("L135").Value = 900012060
$ws.Range("M135").Value = -3511.59375
$ws.Range("N135").Value = -900017130

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12570.833
$ws.Range("I32").Value = 7296.8823
$ws.Range("J32").Value = 22532.74
$ws.Range("K32").Value = 7296.8823
$ws.Range("L32").Value = 22532.74
$ws.Range("M32").Value = -7009.8823
$ws.Range("N32").Value = -23106.74
# Row 45
$ws.Range("H45").Value = 6062268
$ws.Range("I45").Value = 6495158.5
$ws.Range("K45").Value = 6495158.5
$ws.Range("M45").Value = -6494781.5
# Row 74
$ws.Range("H74").Value = 16203.125
$ws.Range("I74").Value = 3100
$ws.Range("J74").Value = 18075
$ws.Range("K74").Value = 3100
$ws.Range("L74").Value = 18075
$ws.Range("M74").Value = -2226
$ws.Range("N74").Value = -19823
# Row 77
$ws.Range("H77").Value = 16203.125
$ws.Range("I77").Value = 3100
$ws.Range("J77").Value = 18075
$ws.Range("K77").Value = 15500
$ws.Range("L77").Value = 90375
$ws.Range("M77").Value = -11132
$ws.Range("N77").Value = -99111
# Row 110
$ws.Range("H110").Value = 838.7727
$ws.Range("I110").Value = 786.3158
$ws.Range("J110").Value = 1171
$ws.Range("K110").Value = 786.3158
$ws.Range("L110").Value = 1171
$ws.Range("M110").Value = 1258.6842
$ws.Range("N110").Value = -5261
# Row 122
$ws.Range("H122").Value = 1897.2727
$ws.Range("I122").Value = 1363.92
$ws.Range("J122").Value = 3564
$ws.Range("K122").Value = 4091.76
$ws.Range("L122").Value = 10692
$ws.Range("M122").Value = -1641.76
$ws.Range("N122").Value = -15592
# Row 132
$ws.Range("H132").Value = 1493.68
$ws.Range("I132").Value = 1130.1794
$ws.Range("J132").Value = 2782.4546
$ws.Range("K132").Value = 3390.5382
$ws.Range("L132").Value = 8347.363799999999
$ws.Range("M132").Value = -860.5382
$ws.Range("N132").Value = -13407.3638

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1234.25
$ws.Range("I107").Value = 1087.3158
$ws.Range("J107").Value = 1544.4445
$ws.Range("K107").Value = 1087.3158
$ws.Range("L107").Value = 1544.4445
$ws.Range("M107").Value = 832.6841999999999
$ws.Range("N107").Value = -5384.4445

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2699.805
$ws.Range("I31").Value = 1756.4
$ws.Range("J31").Value = 3244.077
$ws.Range("K31").Value = 1756.4
$ws.Range("L31").Value = 3244.077
$ws.Range("M31").Value = -1461.4
$ws.Range("N31").Value = -3834.077
# Row 34
$ws.Range("H34").Value = 2699.805
$ws.Range("I34").Value = 1756.4
$ws.Range("J34").Value = 3244.077
$ws.Range("K34").Value = 1756.4
$ws.Range("L34").Value = 3244.077
$ws.Range("M34").Value = -1554.4
$ws.Range("N34").Value = -3648.077
# Row 58
$ws.Range("H58").Value = 4669.2617
$ws.Range("I58").Value = 5087.1924
$ws.Range("J58").Value = 3990.125
$ws.Range("K58").Value = 5087.1924
$ws.Range("L58").Value = 3990.125
$ws.Range("M58").Value = -4884.1924
$ws.Range("N58").Value = -4396.125
# Row 132
$ws.Range("H132").Value = 1936.0212
$ws.Range("I132").Value = 1705.5807
$ws.Range("J132").Value = 2382.5
$ws.Range("K132").Value = 5116.742099999999
$ws.Range("L132").Value = 7147.5
$ws.Range("M132").Value = -2586.742099999999
$ws.Range("N132").Value = -12207.5
# Row 134
$ws.Range("H134").Value = 2598.8438
$ws.Range("I134").Value = 2252.4614
$ws.Range("J134").Value = 4099.8335
$ws.Range("K134").Value = 6757.3842
$ws.Range("L134").Value = 12299.5005
$ws.Range("M134").Value = -4222.3842
$ws.Range("N134").Value = -17369.5005
# Row 136
$ws.Range("H136").Value = 4669.2617
$ws.Range("I136").Value = 5087.1924
$ws.Range("J136").Value = 3990.125
$ws.Range("K136").Value = 15261.5772
$ws.Range("L136").Value = 11970.375
$ws.Range("M136").Value = -12711.5772
$ws.Range("N136").Value = -17070.375

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1122962.6
$ws.Range("J113").Value = 621.15
$ws.Range("L113").Value = 1863.45
$ws.Range("N113").Value = -6203.45

$ws = $wb.Worksheets.Item("GSM")
# Row 63
$ws.Range("H63").Value = 23150.75
$ws.Range("I63").Value = 25103
$ws.Range("J63").Value = 22500
$ws.Range("K63").Value = 25103
$ws.Range("L63").Value = 22500
$ws.Range("M63").Value = -24417
$ws.Range("N63").Value = -23872
# Row 66
$ws.Range("H66").Value = 23150.75
$ws.Range("I66").Value = 25103
$ws.Range("J66").Value = 22500
$ws.Range("K66").Value = 75309
$ws.Range("L66").Value = 67500
$ws.Range("M66").Value = -71877
$ws.Range("N66").Value = -74364
# Row 122
$ws.Range("H122").Value = 869.5
$ws.Range("I122").Value = 869.5
$ws.Range("K122").Value = 2608.5
$ws.Range("M122").Value = -158.5
# Row 126
$ws.Range("H126").Value = 111112520
$ws.Range("J126").Value = 3014
$ws.Range("L126").Value = 9042
$ws.Range("N126").Value = -13982
# Row 132
$ws.Range("H132").Value = 4786.4653
$ws.Range("I132").Value = 3040.2354
$ws.Range("J132").Value = 5928.231
$ws.Range("K132").Value = 9120.706200000001
$ws.Range("L132").Value = 17784.693
$ws.Range("M132").Value = -6590.706200000001
$ws.Range("N132").Value = -22844.693

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2920
$ws.Range("I61").Value = 2046.1538
$ws.Range("K61").Value = 2046.1538
$ws.Range("M61").Value = -1844.1538
# Row 113
$ws.Range("H113").Value = 2920
$ws.Range("I113").Value = 2046.1538
$ws.Range("K113").Value = 2046.1538
$ws.Range("M113").Value = 123.8462
# Row 132
$ws.Range("H132").Value = 3479.0232
$ws.Range("I132").Value = 3356.84
$ws.Range("J132").Value = 3648.7222
$ws.Range("K132").Value = 10070.52
$ws.Range("L132").Value = 10946.1666
$ws.Range("M132").Value = -7540.52
$ws.Range("N132").Value = -16006.1666
# Row 133
$ws.Range("H133").Value = 28550.666
$ws.Range("J133").Value = 28550.666
$ws.Range("L133").Value = 28550.666
$ws.Range("N133").Value = -33610.666

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 1580
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 1966.6666
$ws.Range("K62").Value = 1000
$ws.Range("L62").Value = 1966.6666
$ws.Range("M62").Value = -376
$ws.Range("N62").Value = -3214.6666
# Row 65
$ws.Range("H65").Value = 1580
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 1966.6666
$ws.Range("K65").Value = 5000
$ws.Range("L65").Value = 9833.333000000001
$ws.Range("M65").Value = -1880
$ws.Range("N65").Value = -16073.333
# Row 132
$ws.Range("H132").Value = 935850.9399999999
$ws.Range("I132").Value = 1402324
$ws.Range("J132").Value = 2904.9048
$ws.Range("K132").Value = 4206972
$ws.Range("L132").Value = 8714.714399999999
$ws.Range("M132").Value = -4204442
$ws.Range("N132").Value = -13774.7144
